$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Ordered list of (old, new) expression pairs, matching the table cells
# in row-major order (20 rows x 5 columns = 100 cells).
$pairs = @(
    @("32-22=10", "4+14=18"),
    @("5+50=55", "9+37=46"),
    @("6+13=19", "82-9=73"),
    @("30-25=5", "89-2=87"),
    @("81-4=77", "88+9=97"),
    @("18+35=53", "95-5=90"),
    @("69-13=56", "14-13=1"),
    @("71-19=52", "37-32=5"),
    @("16+64=80", "15-7=8"),
    @("50-43=7", "79+14=93"),
    @("39+16=55", "84-59=25"),
    @("38-28=10", "10+13=23"),
    @("33-3=30", "13+50=63"),
    @("59+11=70", "36+38=74"),
    @("37+8=45", "20+51=71"),
    @("38+2=40", "50+35=85"),
    @("80-9=71", "56-10=46"),
    @("22+10=32", "57+0=57"),
    @("76-72=4", "69-0=69"),
    @("79-25=54", "8+50=58"),
    @("71-4=67", "11+54=65"),
    @("44+14=58", "3+13=16"),
    @("5+79=84", "94-87=7"),
    @("36+13=49", "19-19=0"),
    @("0+1=1", "46-10=36"),
    @("24+24=48", "2+60=62"),
    @("11-9=2", "55+38=93"),
    @("84-28=56", "32-16=16"),
    @("87-86=1", "25-18=7"),
    @("47+48=95", "57-55=2"),
    @("82-25=57", "3+19=22"),
    @("38-13=25", "21+33=54"),
    @("61+4=65", "40+18=58"),
    @("46-11=35", "74-57=17"),
    @("27+34=61", "46-36=10"),
    @("52-21=31", "42-30=12"),
    @("57+33=90", "78-41=37"),
    @("33+53=86", "79-23=56"),
    @("12+32=44", "10+59=69"),
    @("40-9=31", "74+8=82"),
    @("59-51=8", "75-32=43"),
    @("70-14=56", "80+9=89"),
    @("23+50=73", "51+40=91"),
    @("20+66=86", "11+69=80"),
    @("39+35=74", "2+45=47"),
    @("7+25=32", "20-4=16"),
    @("87+8=95", "36-30=6"),
    @("47+7=54", "45+42=87"),
    @("33+18=51", "54-49=5"),
    @("82-35=47", "78-20=58"),
    @("31+46=77", "32-18=14"),
    @("98-54=44", "64-25=39"),
    @("92-54=38", "24+66=90"),
    @("25+6=31", "1+34=35"),
    @("66+15=81", "93-77=16"),
    @("42-32=10", "76-27=49"),
    @("28+63=91", "8+83=91"),
    @("97-36=61", "20+22=42"),
    @("24+2=26", "95-39=56"),
    @("62+28=90", "48-41=7"),
    @("98-11=87", "49+46=95"),
    @("75+13=88", "37-22=15"),
    @("68-3=65", "79-53=26"),
    @("12+83=95", "61-27=34"),
    @("9+47=56", "85-63=22"),
    @("98-51=47", "86-27=59"),
    @("9+56=65", "91-47=44"),
    @("20+75=95", "6+16=22"),
    @("85-48=37", "72-17=55"),
    @("75+0=75", "95-52=43"),
    @("15-5=10", "37-12=25"),
    @("63+9=72", "72-8=64"),
    @("46+44=90", "76+6=82"),
    @("77-71=6", "81-1=80"),
    @("85+14=99", "50-27=23"),
    @("4+78=82", "60+20=80"),
    @("9-3=6", "51+33=84"),
    @("43+52=95", "6+86=92"),
    @("95-9=86", "29+66=95"),
    @("36-36=0", "16+50=66"),
    @("5+74=79", "96-41=55"),
    @("66+5=71", "86-56=30"),
    @("68-0=68", "99-34=65"),
    @("30-13=17", "83-62=21"),
    @("19+48=67", "57-13=44"),
    @("37+50=87", "84-4=80"),
    @("23+36=59", "26+17=43"),
    @("9+19=28", "25+11=36"),
    @("52-21=31", "81-29=52"),
    @("40+15=55", "84-47=37"),
    @("24-20=4", "19+40=59"),
    @("52+45=97", "73+4=77"),
    @("41+25=66", "28+60=88"),
    @("14+19=33", "83+8=91"),
    @("86-78=8", "76+11=87"),
    @("3+34=37", "81-33=48"),
    @("10-3=7", "40-33=7"),
    @("3+28=31", "6+5=11"),
    @("92+1=93", "10+3=13"),
    @("14+1=15", "76-7=69"),
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
$mismatches = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $expected = $pairs[$idx][0]
        $replacement = $pairs[$idx][1]
        $cell = $t.Cell($r, $c)
        $current = $cell.Range.Text
        if ($current.Substring(0, $expected.Length) -ne $expected) {
            $mismatches++
            Write-Host "Mismatch at row $r col $c : expected '$expected' got '$current'"
        }
        $cell.Range.Text = $replacement
        $idx++
    }
}

Write-Host "Updated $idx cells. Mismatches: $mismatches"
